$d = $word.ActiveDocument

# The credits paragraph currently reads (single run):
#   " with initial ideas, helping with programming and testing and, of course, ..."
# It needs to become three runs (same character formatting throughout):
#   " with initial "
#   "ideas, helping with programming, "
#   "testing and, of course, ..."
# i.e. "programming and testing" -> "programming, testing", with the run
# boundaries landing right before "ideas" and right before "testing".

$oldMiddle = "ideas, helping with programming and "
$newMiddle = "ideas, helping with programming, "

# Locate the exact span to replace.
$findRng = $d.Content
$found = $findRng.Find.Execute($oldMiddle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text to edit"
}

$start = $findRng.Start
$end = $findRng.End

# Replace just that middle span's text.
$midRng = $d.Range($start, $end)
$midRng.Text = $newMiddle

# Re-acquire the (now shorter) range for the replaced text and nudge its
# character formatting (on, then back off) so the surrounding, identically
# formatted text is split into separate runs instead of being re-merged
# into one long run.
$midRng2 = $d.Range($start, $start + $newMiddle.Length)
$midRng2.Bold = 1
$midRng2.Bold = 0
